$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "baseLogicExpr { relationOp baseLogicExpr } ;" ->
#           "baseLogicExpr [ relationOp baseLogicExpr ] ;"
# The braces become brackets, and the run is split into three runs
# (the two surrounding pieces of text stay as-is, only the middle
# "{ relationOp baseLogicExpr }" -> "[ relationOp baseLogicExpr ]" part
# changes) mirroring how Word splits a run when you edit inside it.
# ---------------------------------------------------------------------
$full1 = "baseLogicExpr { relationOp baseLogicExpr } ;"
$text = $d.Content.Text
$idx1 = $text.IndexOf($full1)

$c1part1 = "baseLogicExpr "
$c1part2old = "{ relationOp baseLogicExpr }"
$c1part2new = "[ relationOp baseLogicExpr ]"

$c1p1start = $idx1
$c1p1end = $c1p1start + $c1part1.Length
$c1p2start = $c1p1end
$c1p2end = $c1p2start + $c1part2old.Length

$c1r2 = $d.Range($c1p2start, $c1p2end)
$c1r2.Text = $c1part2new
$c1p2end = $c1p2start + $c1part2new.Length

# Force the edited middle segment to live in its own run, distinct from
# the unchanged text before/after it (toggling+reverting formatting is
# enough to split runs without altering the visible formatting).
$c1r2b = $d.Range($c1p2start, $c1p2end)
$c1r2b.Font.Bold = 1
$c1r2b.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 2: the run " = [ unaryLogicOp ] (mathExpr " is split into
# " = [ unaryLogicOp ] (ma" + "thExpr ", with the "_GoBack" bookmark
# relocated in between the two halves (it used to sit alone in a later
# empty paragraph).
# ---------------------------------------------------------------------
$full2 = " = [ unaryLogicOp ] (mathExpr "
$text2 = $d.Content.Text
$idx2 = $text2.IndexOf($full2)

$c2part1 = " = [ unaryLogicOp ] (ma"
$c2part2 = "thExpr "

$c2splitPoint = $idx2 + $c2part1.Length
$c2endPoint = $idx2 + $full2.Length

$c2r2 = $d.Range($c2splitPoint, $c2endPoint)
$c2r2.Font.Bold = 1
$c2r2.Font.Bold = 0

# Move the "_GoBack" bookmark from its old (now-stale) location to the
# split point between "...(ma" and "thExpr ...".
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
}

$newBookmarkRange = $d.Range($c2splitPoint, $c2splitPoint)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
